$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Delete the two vestigial "loop template" rows (old rows 7 and 8); this
# shifts the real response rows (old 9:17) up to become rows 7:15 and
# updates the used-range dimension accordingly.
$ws.Rows("7:8").Delete()

# Row 5: a new blank (bold-styled) cell appears at AW5.
$ws.Range("AW5").Font.Bold = $true

# Row 6: the "0" that was in AV6 moves one cell to the left, into AS6.
$ws.Range("AS6").Value = 0
$ws.Range("AV6").ClearContents()

# The hyperlink that used to sit on U11 now lives on U9 (its row shifted
# up by two along with the rest of that response row). Recreate it there.
$ws.Range("U11").Hyperlinks.Delete()
$hlink = $ws.Hyperlinks.Add($ws.Range("U9"), "http://coyoteoutdoorschool.org")
$hlink.Address = "http://coyoteoutdoorschool.org"
$ws.Range("U9").Style = "Hyperlink"
